$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 50,4

$arr[0,0] = "Bitcoin"
$arr[0,1] = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$arr[0,2] = "30.356.75"
$arr[0,3] = "  +0.41%  "

$arr[1,0] = "Ethereum"
$arr[1,1] = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$arr[1,2] = "1.870.36"
$arr[1,3] = "  -0.30%  "

$arr[2,0] = "TetherUSD"
$arr[2,1] = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$arr[2,2] = "'1.001"
$arr[2,3] = "  -0.03%  "

$arr[3,0] = "BNB"
$arr[3,1] = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$arr[3,2] = "'237.93"
$arr[3,3] = "  +0.99%  "

$arr[4,0] = "USDC"
$arr[4,1] = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$arr[4,2] = "'1.001"
$arr[4,3] = "  +0.00%  "

$arr[5,0] = "XRP"
$arr[5,1] = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$arr[5,2] = "'0.4813"
$arr[5,3] = "  -0.31%  "

$arr[6,0] = "Cardano"
$arr[6,1] = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$arr[6,2] = "'0.2809"
$arr[6,3] = "  -2.13%  "

$arr[7,0] = "Dogecoin"
$arr[7,1] = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$arr[7,2] = "'0.06508"
$arr[7,3] = "  -0.89%  "

$arr[8,0] = "WrappedEther"
$arr[8,1] = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$arr[8,2] = "1.872.08"
$arr[8,3] = "  -0.27%  "

$arr[9,0] = "TRON"
$arr[9,1] = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$arr[9,2] = "'0.07482"
$arr[9,3] = "  +2.33%  "

$arr[10,0] = "Solana"
$arr[10,1] = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$arr[10,2] = "'16.50"
$arr[10,3] = "  -0.99%  "

$arr[11,0] = "Polkadot"
$arr[11,1] = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$arr[11,2] = "'5.086"
$arr[11,3] = "  -0.93%  "

$arr[12,0] = "Litecoin"
$arr[12,1] = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$arr[12,2] = "'87.92"
$arr[12,3] = "  +1.35%  "

$arr[13,0] = "Polygon"
$arr[13,1] = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$arr[13,2] = "'0.6559"
$arr[13,3] = "  +0.65%  "

$arr[14,0] = "WrappedBTC"
$arr[14,1] = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$arr[14,2] = "30.345.23"
$arr[14,3] = "  +0.47%  "

$arr[15,0] = "Avalanche"
$arr[15,1] = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$arr[15,2] = "'13.25"
$arr[15,3] = "  -0.34%  "

$arr[16,0] = "Dai"
$arr[16,1] = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$arr[16,2] = "'1.001"
$arr[16,3] = "  +0.00%  "

$arr[17,0] = "ShibaInu"
$arr[17,1] = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$arr[17,2] = "'0.000007612"
$arr[17,3] = "  -1.33%  "

$arr[18,0] = "WrappedliquidstakedEther2.0"
$arr[18,1] = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$arr[18,2] = "2.116.60"
$arr[18,3] = "  -0.32%  "

$arr[19,0] = "BitcoinCash"
$arr[19,1] = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$arr[19,2] = "'222.67"
$arr[19,3] = "  +15.88%  "

$arr[20,0] = "BinanceUSD"
$arr[20,1] = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$arr[20,2] = "'1.001"
$arr[20,3] = "  -0.11%  "

$arr[21,0] = "Uniswap"
$arr[21,1] = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$arr[21,2] = "'5.288"
$arr[21,3] = "  -0.38%  "

$arr[22,0] = "Chainlink"
$arr[22,1] = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$arr[22,2] = "'6.167"
$arr[22,3] = "  +1.21%  "

$arr[23,0] = "Cosmos"
$arr[23,1] = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$arr[23,2] = "'9.277"
$arr[23,3] = "  +0.33%  "

$arr[24,0] = "Monero"
$arr[24,1] = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$arr[24,2] = "'167.04"
$arr[24,3] = "  +3.33%  "

$arr[25,0] = "EthereumClassic"
$arr[25,1] = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$arr[25,2] = "'18.63"
$arr[25,3] = "  +3.88%  "

$arr[26,0] = "LidoDAOToken"
$arr[26,1] = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$arr[26,2] = "'1.985"
$arr[26,3] = "  +4.37%  "

$arr[27,0] = "Toncoin"
$arr[27,1] = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$arr[27,2] = "'1.461"
$arr[27,3] = "  +1.68%  "

$arr[28,0] = "Stellar"
$arr[28,1] = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$arr[28,2] = "'0.09350"
$arr[28,3] = "  +3.03%  "

$arr[29,0] = "InternetComputer(DFINITY)"
$arr[29,1] = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$arr[29,2] = "'4.304"
$arr[29,3] = "  +1.33%  "

$arr[30,0] = "Filecoin"
$arr[30,1] = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$arr[30,2] = "'4.014"
$arr[30,3] = "  +0.43%  "

$arr[31,0] = "Hedera"
$arr[31,1] = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$arr[31,2] = "'0.05039"
$arr[31,3] = "  -0.41%  "

$arr[32,0] = "ARBITRUM"
$arr[32,1] = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$arr[32,2] = "'1.211"
$arr[32,3] = "  +10.99%  "

$arr[33,0] = "ImmutableX"
$arr[33,1] = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$arr[33,2] = "'0.7491"
$arr[33,3] = "  +5.32%  "

$arr[34,0] = "HuobiToken"
$arr[34,1] = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$arr[34,2] = "'2.711"
$arr[34,3] = "  +0.36%  "

$arr[35,0] = "VeChain"
$arr[35,1] = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$arr[35,2] = "'0.01832"
$arr[35,3] = "  +2.83%  "

$arr[36,0] = "MXToken"
$arr[36,1] = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$arr[36,2] = "'2.620"
$arr[36,3] = "  -0.46%  "

$arr[37,0] = "RenderToken"
$arr[37,1] = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$arr[37,2] = "'2.082"
$arr[37,3] = "  +2.35%  "

$arr[38,0] = "TrustWalletToken"
$arr[38,1] = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$arr[38,2] = "'0.9048"
$arr[38,3] = "  -1.51%  "

$arr[39,0] = "Quant"
$arr[39,1] = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$arr[39,2] = "'106.73"
$arr[39,3] = "  +1.05%  "

$arr[40,0] = "FraxShare"
$arr[40,1] = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$arr[40,2] = "'5.918"
$arr[40,3] = "  +2.49%  "

$arr[41,0] = "TheSandbox"
$arr[41,1] = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$arr[41,2] = "'0.4277"
$arr[41,3] = "  +0.89%  "

$arr[42,0] = "PaxDollar"
$arr[42,1] = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$arr[42,2] = "'1.003"
$arr[42,3] = "  +0.30%  "

$arr[43,0] = "Aptos"
$arr[43,1] = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$arr[43,2] = "'7.417"
$arr[43,3] = "  +0.86%  "

$arr[44,0] = "Algorand"
$arr[44,1] = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$arr[44,2] = "'0.1296"
$arr[44,3] = "  -1.07%  "

$arr[45,0] = "Aave"
$arr[45,1] = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$arr[45,2] = "'64.07"
$arr[45,3] = "  -1.12%  "

$arr[46,0] = "NEARProtocol"
$arr[46,1] = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$arr[46,2] = "'1.482"
$arr[46,3] = "  +9.25%  "

$arr[47,0] = "EnergySwap"
$arr[47,1] = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$arr[47,2] = "'8.992"
$arr[47,3] = "  +1.17%  "

$arr[48,0] = "Elrond"
$arr[48,1] = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$arr[48,2] = "'33.89"
$arr[48,3] = "  +0.74%  "

$arr[49,0] = "Cronos"
$arr[49,1] = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$arr[49,2] = "'0.05655"
$arr[49,3] = "  -1.55%  "

$ws.Range("B2:E51").Value = $arr